$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.055.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.418.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.110'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.24'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000174'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.848.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.950.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.412.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '321.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '565.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.08%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0935'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '148.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0532'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.593'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0920'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E51").Value = '  +0.59%  '
